$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 280 (shifts existing row 280 onward down to 281+)
$ws.Rows(280).Insert()

# Populate the newly inserted row 280 with the new data record
$ws.Cells.Item(280, 1).Value = 3
$ws.Cells.Item(280, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(280, 3).Value = "Coquimbo"
$ws.Cells.Item(280, 4).Value = 44559
$ws.Cells.Item(280, 5).Value = 5
$ws.Cells.Item(280, 6).Value = 100112037
$ws.Cells.Item(280, 7).Value = "Cebollín"
$ws.Cells.Item(280, 8).Value = "Sin especificar"
$ws.Cells.Item(280, 9).Value = "Primera"
$ws.Cells.Item(280, 10).Value = 198
$ws.Cells.Item(280, 11).Value = 3000
$ws.Cells.Item(280, 12).Value = 3500
$ws.Cells.Item(280, 13).Value = 3247
$ws.Cells.Item(280, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(280, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(280, 16).Value = 90
$ws.Cells.Item(280, 17).Value = 36
$ws.Cells.Item(280, 18).Value = "Hortaliza"
